# Add new Upgrade cards to the CardDB (3rd worksheet) of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

$rows = @(
    @('U0071', 'Collision Detector', 'When performing a boost, barrel roll, or decloak, your ship and maneuver template can overlap obstacles. When rolling for obstacle damage, ignore all CRIT results.', 0, 'N', 'N', 'System', 1),
    @('U0072', 'Sensor Cluster', 'When defending, you may spend a focus token to change one of your blank results to an EVADE result.', 2, 'N', 'N', 'Tech', 2),
    @('U0073', 'Special Ops Training', 'TIE/SF ONLY. When attacking with a primary weapon from your primary firing arc, you may roll 1 additional attack die. If you do not, you may perform an additional attack from your auxiliary firing arc.', 0, 'N', 'N', 'Title', 1),
    @('U0074', 'Bodyguard', 'SCUM ONLY. At the start of the combat phase, you may spend a focus token to choose a friendly ship at Range 1 with higher pilot skill than you. Increase its agility value by 1 until the end of the round.', 2, 'Y', 'N', 'Elite Pilot Skill', 1),
    @('U0075', 'Calculation', 'When attacking, you may spend 1 focus token to change 1 of your FOCUS results to a CRIT result.', 1, 'N', 'N', 'Elite Pilot Skill', 1),
    @('U0076', 'Accuracy Corrector', 'When attacking, during the ''Modify Attack Dice'' step, you may cancel all your dice results. If you do, you may add 2 HIT results to your roll. Your dice cannot be modified again during this attack.', 3, 'N', 'N', 'System', 1),
    @('U0077', 'Ion Torpedoes', 'FP: 4, RNG: 2-3 ATTACK (TARGET LOCK): Spend your target lock and discard this card to perform this attack. If this attack hits, the defender and each ship at range 1 of it receive 1 ion token.', 5, 'N', 'N', 'Torpedo', 3),
    @('U0078', 'Inertial Dampeners', 'When you reveal your maneuver, you may discard this card to instead perform a white [STAY 0] maneuver. Then receive 1 stress token.', 1, 'N', 'N', 'Illicit', 1),
    @('U0079', 'Autothrusters', 'When defending, if you are beyond range 2 or outside the attacker''s firing arc, you may change 1 of your blank results to an EVADE result. You can equip this card only if you have the BOOST action icon.', 2, 'N', 'N', 'Modification', 2),
    @('U0080', 'Hull Upgrade', 'Increase your hull value by 1.', 3, 'N', 'N', 'Modification', 1),
    @('U0081', 'Virago', 'STARVIPER ONLY. Your upgrade bar gains the SYSTEM and ILLICIT upgrade icons. You cannot equip this card if your pilot skill value is ''3'' or lower.', 1, 'Y', 'N', 'Title', 1),
    @('U0082', 'Homing Missiles', 'FP: 4, RNG: 2-3 ATTACK (TARGET LOCK): Discard this card to perform this attack. The defender cannot spend evade tokens during this attack.', 5, 'N', 'N', 'Missile', 1),
    @('U0083', 'Assault Missiles', 'FP: 4, RNG: 2-3 ATTACK (TARGET LOCK): Spend your target lock and discard this card to perform this attack. If this attack hits, each other ship at Range 1 of the defender suffers 1 damage.', 5, 'N', 'N', 'Missile', 1),
    @('U0084', 'Expose', 'ACTION: Until the end of the round, increase your primary weapon value 1, and decrease your agility value by 1.', 4, 'N', 'N', 'Elite Pilot Skill', 1),
    @('U0085', 'Veteran Instincts', 'Increase your pilot skill value by 2.', 1, 'N', 'N', 'Elite Pilot Skill', 1),
    @('U0086', 'Seismic Charges', 'When you reveal you maneuver dial, you may discard this card to DROP 1 seismic charge token. This token DETONATES at the end of the activation phase.', 2, 'N', 'N', 'Bomb', 1),
    @('U0087', 'Stealth Device', 'Increase your agility value by 1. If you are hit by an attack, discard this card.', 3, 'N', 'N', 'Modification', 2),
    @('U0088', 'Mercenary Copilot', 'When attacking at Range 3, you may change 1 of your HIT results to a CRIT result.', 2, 'N', 'N', 'Crew', 1),
    @('U0089', 'Gunner', 'After you perform and attack that does not hit, you may immediately perform a primary weapon attack. You cannot perform another attack this round.', 5, 'N', 'N', 'Crew', 1),
    @('U0090', 'Slave 1', 'FIRESPRAY-31 ONLY. Your upgrade bar gains the TORPEDO upgrade icon.', 0, 'Y', 'N', 'Title', 1),
    @('U0091', 'Ruthlessness', 'IMPERIAL ONLY. After you perform an attack that hits, you MUST choose 1 other ship at Range 1 of the defender (other than yourself). That ship suffers 1 damage.', 3, 'N', 'N', 'Elite Pilot Skill', 2),
    @('U0092', 'Intimidation', 'While you are touching an enemy ship, reduce that ship''s agility value by 1.', 2, 'N', 'N', 'Elite Pilot Skill', 1),
    @('U0093', 'Fleet Officer', 'IMPERIAL ONLY. ACTION: choose up to 2 friendly ships at Range 1-2 and assign 1 focus token to each of those ships. Then receive 1 stress token.', 3, 'N', 'N', 'Crew', 1),
    @('U0094', 'Mara Jade', 'IMPERIAL ONLY. At the end of the Combat phase, each enemy shipat Range 1 that does not have a stress token receives 1 stress token.', 3, 'Y', 'N', 'Crew', 1),
    @('U0095', 'Ysanne Isard', 'IMPERIAL ONLY. At the start of the Combat phase, if you have no shields and at least 1 Damage card assigned to your ship, you may perform a free evade action.', 4, 'Y', 'N', 'Crew', 1),
    @('U0096', 'Moff Jerjerrod', 'IMPERIAL ONLY. When you are dealt a faceup Damage card, you may discard this Upgrade card or another CREW Upgrade card to flip that Damage card facedown (without resolving its effect).', 2, 'Y', 'N', 'Crew', 1),
    @('U0097', 'Proton Bombs', 'When you reveal your maneuver dial, you may discard this card to DROP 1 proton bomb token. This token DETONATES at the end of the Activation phase.', 5, 'N', 'N', 'Bomb', 1),
    @('U0098', 'Tactical Jammer', 'LARGE SHIP ONLY. Your ship can obstruct enemy attacks.', 1, 'N', 'N', 'Modification', 2),
    @('U0099', 'Dauntless', 'VT-49 DECIMATOR ONLY. After you execute a maneuver that causes you to overlap another ship, you may perform 1 free action. Then receive 1 stress token.', 2, 'Y', 'N', 'Title', 1),
    @('', 'Trick Shot', 'When attacking, if the attack is obstructed, you may roll an additional attack die.', 0, 'N', 'N', 'Elite Pilot Skill', 2),
    @('', 'Finn', 'REBEL ONLY. When attacking with a primary weapon or defending, if the enemy ship is inside your firing arc, you may add 1 blank result to your roll.', 5, 'Y', 'N', 'Crew', 1),
    @('', 'Rey', 'REBEL ONLY. At the start of the End phase, you may place one of your ship''s focus tokens on this card. At the start of the Combat phase, you may assign 1 of those tokens to your ship.', 2, 'Y', 'N', 'Crew', 1),
    @('', 'Hotshot Co-Pilot', 'When attacking with a primary weapon, the defender must spend 1 focus token if able. When defending, the attacker must spend 1 focus token if able.', 4, 'N', 'N', 'Crew', 1),
    @('', 'Snap Shot', 'FP: 2, RNG: 1  After an enemy ship executes a maneuver, you may perform this attack against that ship. ATTACK: Attack 1 ship. You cannot modify your attack dice and cannot attack again this phase.', 2, 'N', 'N', 'Elite Pilot Skill', 2),
    @('', 'M9-G8', 'When a ship you have locked is attacking, you may choose 1 attack die. The attacker must reroll that die. You can acquire target locks on other friendly ships.', 3, 'Y', 'N', 'Astromech', 1),
    @('', 'Burnout SLAM', 'LARGE SHIP ONLY. Your action bar gains the SLAM action icon. After you perform a SLAM action, discard this card.', 1, 'N', 'N', 'Illicit', 2),
    @('', 'Primed Thrusters', 'SMALL SHIP ONLY. Stress tokens do not prevent you from performing boost or barrel roll actions unless you have 3 or more stress tokens.', 1, 'N', 'N', 'Tech', 1),
    @('', 'Pattern Analyzer', 'When executing a maneuver, you may resolve the ''Check Pilot Stress'' step after the ''Perform Action'' step (instead of before that step).', 2, 'N', 'N', 'Tech', 2),
    @('', 'Millenium Falcon', 'YT-1300 ONLY. After you execute a 3-speed bank maneuver, if you are not touching another ship and you are not stressed, you may receive 1 stress token to rotate your ship 180 degrees.', 1, 'Y', 'N', 'Title', 1),
    @('', 'Black One', 'T-70 X-WING ONLY. After you perform a boost or barrel roll action, you may remove 1 enemy target lock from a friendly ship at range 1. You cannot equip this card if your pilot skill is ''6'' or lower.', 1, 'Y', 'N', 'Title', 1),
    @('', 'Smuggling Compartment', 'YT-1300 OR YT-2400 ONLY. Your upgrade bar gains the ILLICIT icon. You may equip 1 additional Modification upgrade that costs 3 or fewer squad points.', 0, 'N', 'Y', 'Modification', 1)
)

$startRow = 72

# Copy formatting (styles) from the last existing data row (row 71) down across
# the whole new block so the new rows look like the rest of the table
# (vertical-center on most columns, wrap-text on the Text column).
$lastRow = $startRow + $rows.Length - 1
$ws.Range("A71:H71").Copy($ws.Range("A" + $startRow + ":H" + $lastRow))

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $row = $startRow + $i

    $id = $r[0]
    if ($id -ne "") {
        $ws.Cells.Item($row, 1).Value = $id
    } else {
        $ws.Cells.Item($row, 1).Value = ""
    }

    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
}

# A handful of pre-existing rows also had their Qty (column H) bumped.
$ws.Cells.Item(41, 8).Value = 3
$ws.Cells.Item(49, 8).Value = 2
$ws.Cells.Item(60, 8).Value = 2
$ws.Cells.Item(61, 8).Value = 2
$ws.Cells.Item(65, 8).Value = 2
